$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Prices in column D are plain text (e.g. "46.784.90") so they must be
# written with a temporary text NumberFormat to stop Excel auto-converting
# them into numbers (which would also strip trailing zeros). The style is
# reset back to Normal right after so no visible formatting changes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.784.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.337.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.24%  "

$ws.Range("E4").Value = "  -0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("E12").Value = "  +3.11%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.693.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.342.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.24%  "

$ws.Range("E16").Value = "  +4.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.828"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.697.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +17.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.49%  "

$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0816"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.29%  "

$ws.Range("E40").Value = "  +5.20%  "

$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.80%  "

$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.841.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.07%  "

$ws.Range("E46").Value = "  +5.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.98%  "
